# Powerpoint writer: consolidate text run nodes.
# This should reduce the size of the generated files.
#
# Merges the first two runs of specific paragraphs ("Slide" + " " -> "Slide ",
# and "an" + " " -> "an ") into a single run, leaving the following run
# ("1"/"2"/"image") untouched, matching the target diff exactly.

$p = $ppt.ActivePresentation

function Merge-FirstTwoRuns($shape, $mergedText, $charCount) {
    $tr = $shape.TextFrame.TextRange
    $sub = $tr.Characters(1, $charCount)
    $sub.Text = $mergedText
}

# Slide 1: Title "Slide" + " " + "1" -> "Slide " + "1"
$slide1 = $p.Slides.Item(1)
Merge-FirstTwoRuns $slide1.Shapes.Item(1) "Slide " 6

# Slide 1: TextBox "an" + " " + "image" -> "an " + "image"
Merge-FirstTwoRuns $slide1.Shapes.Item(3) "an " 3

# Slide 2: Title "Slide" + " " + "2" -> "Slide " + "2"
$slide2 = $p.Slides.Item(2)
Merge-FirstTwoRuns $slide2.Shapes.Item(1) "Slide " 6

# Slide 2: TextBox "an" + " " + "image" -> "an " + "image"
Merge-FirstTwoRuns $slide2.Shapes.Item(4) "an " 3
